# "Generate Report for Handback"
#
# The handback transform failed for the 7b5c8606-... entry (row 3 of the
# status tables) because the handback file name didn't match the handoff
# file name. This updates the status from "Ready for handoff" to
# "Handback transform failed" everywhere it is shown (Overview summary
# sheet plus each per-locale detail sheet), and records the mismatch
# reason in the "Error Detail" column (K) of the per-locale sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 (the 7b5c8606-... file) shows the status for both
# the zh-cn and de-de columns (B and C).
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# Per-locale detail sheets: row 3's Status column (C).
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Per-locale detail sheets: populate the new Error Detail (K) cell for
# row 3 explaining the handback/handoff file name mismatch.
$wsZhCn.Range("K3").Value = "Handback file name: bibyn0cr.cqs is different with handoff file name: 7b5c8606-716c-4e52-8cd9-cbefe60558a1.cbdd3b0e02534cb1c1ca6ae40e32cecc8b488ea7.zh-cn."
$wsDeDe.Range("K3").Value = "Handback file name: bibyn0cr.cqs is different with handoff file name: 7b5c8606-716c-4e52-8cd9-cbefe60558a1.cbdd3b0e02534cb1c1ca6ae40e32cecc8b488ea7.de-de."
